$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "13+50="
$t.Cell(1,2).Range.Text = "98-71="
$t.Cell(1,3).Range.Text = "74-58="
$t.Cell(1,4).Range.Text = "94-41="
$t.Cell(1,5).Range.Text = "23+7="

$t.Cell(2,1).Range.Text = "51-44="
$t.Cell(2,2).Range.Text = "75+0="
$t.Cell(2,3).Range.Text = "82-2="
$t.Cell(2,4).Range.Text = "12+3="
$t.Cell(2,5).Range.Text = "64-16="

$t.Cell(3,1).Range.Text = "74-3="
$t.Cell(3,2).Range.Text = "76-72="
$t.Cell(3,3).Range.Text = "31+9="
$t.Cell(3,4).Range.Text = "76-24="
$t.Cell(3,5).Range.Text = "19+7="

$t.Cell(4,1).Range.Text = "77-37="
$t.Cell(4,2).Range.Text = "41+11="
$t.Cell(4,3).Range.Text = "95-74="
$t.Cell(4,4).Range.Text = "67-46="
$t.Cell(4,5).Range.Text = "96-59="

$t.Cell(5,1).Range.Text = "38+48="
$t.Cell(5,2).Range.Text = "80+15="
$t.Cell(5,3).Range.Text = "32+28="
$t.Cell(5,4).Range.Text = "57+7="
$t.Cell(5,5).Range.Text = "47-26="

$t.Cell(6,1).Range.Text = "90-14="
$t.Cell(6,2).Range.Text = "89-73="
$t.Cell(6,3).Range.Text = "15+8="
$t.Cell(6,4).Range.Text = "4+12="
$t.Cell(6,5).Range.Text = "21-9="

$t.Cell(7,1).Range.Text = "91-9="
$t.Cell(7,2).Range.Text = "18+64="
$t.Cell(7,3).Range.Text = "99-18="
$t.Cell(7,4).Range.Text = "11+10="
$t.Cell(7,5).Range.Text = "69+13="

$t.Cell(8,1).Range.Text = "69+7="
$t.Cell(8,2).Range.Text = "72+9="
$t.Cell(8,3).Range.Text = "48-28="
$t.Cell(8,4).Range.Text = "80-35="
$t.Cell(8,5).Range.Text = "61+31="

$t.Cell(9,1).Range.Text = "16-8="
$t.Cell(9,2).Range.Text = "16+81="
$t.Cell(9,3).Range.Text = "91-28="
$t.Cell(9,4).Range.Text = "47-40="
$t.Cell(9,5).Range.Text = "2+27="

$t.Cell(10,1).Range.Text = "11-0="
$t.Cell(10,2).Range.Text = "5+14="
$t.Cell(10,3).Range.Text = "12+18="
$t.Cell(10,4).Range.Text = "29+38="
$t.Cell(10,5).Range.Text = "57-5="

$t.Cell(11,1).Range.Text = "75-40="
$t.Cell(11,2).Range.Text = "98-83="
$t.Cell(11,3).Range.Text = "95-13="
$t.Cell(11,4).Range.Text = "50-30="
$t.Cell(11,5).Range.Text = "85-43="

$t.Cell(12,1).Range.Text = "63-17="
$t.Cell(12,2).Range.Text = "84-23="
$t.Cell(12,3).Range.Text = "27+52="
$t.Cell(12,4).Range.Text = "53+20="
$t.Cell(12,5).Range.Text = "12+44="

$t.Cell(13,1).Range.Text = "16+50="
$t.Cell(13,2).Range.Text = "89-77="
$t.Cell(13,3).Range.Text = "85-40="
$t.Cell(13,4).Range.Text = "56-6="
$t.Cell(13,5).Range.Text = "18+30="

$t.Cell(14,1).Range.Text = "13+19="
$t.Cell(14,2).Range.Text = "70-65="
$t.Cell(14,3).Range.Text = "98-94="
$t.Cell(14,4).Range.Text = "96-61="
$t.Cell(14,5).Range.Text = "9+57="

$t.Cell(15,1).Range.Text = "84-3="
$t.Cell(15,2).Range.Text = "35+24="
$t.Cell(15,3).Range.Text = "83+4="
$t.Cell(15,4).Range.Text = "37-4="
$t.Cell(15,5).Range.Text = "91-59="

$t.Cell(16,1).Range.Text = "92-53="
$t.Cell(16,2).Range.Text = "34-5="
$t.Cell(16,3).Range.Text = "87-19="
$t.Cell(16,4).Range.Text = "33+9="
$t.Cell(16,5).Range.Text = "28-2="

$t.Cell(17,1).Range.Text = "83+10="
$t.Cell(17,2).Range.Text = "68+8="
$t.Cell(17,3).Range.Text = "5+81="
$t.Cell(17,4).Range.Text = "74+18="
$t.Cell(17,5).Range.Text = "90-43="

$t.Cell(18,1).Range.Text = "81-0="
$t.Cell(18,2).Range.Text = "73-32="
$t.Cell(18,3).Range.Text = "86-43="
$t.Cell(18,4).Range.Text = "0+81="
$t.Cell(18,5).Range.Text = "67-40="

$t.Cell(19,1).Range.Text = "69-58="
$t.Cell(19,2).Range.Text = "79-58="
$t.Cell(19,3).Range.Text = "82-60="
$t.Cell(19,4).Range.Text = "13+11="
$t.Cell(19,5).Range.Text = "14+78="

$t.Cell(20,1).Range.Text = "11+7="
$t.Cell(20,2).Range.Text = "51-15="
$t.Cell(20,3).Range.Text = "58-24="
$t.Cell(20,4).Range.Text = "78-73="
$t.Cell(20,5).Range.Text = "45-22="

